$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.423.47"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "'1.667.01"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'313.43"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.3962"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").Value = "'0.3936"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "'52.29"
$ws.Range("E9").Value = "  +6.92%  "
$ws.Range("D10").Value = "'1.396"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'0.08585"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("D14").Value = "'7.333"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "'7.940"
$ws.Range("E15").Value = "  +6.36%  "
$ws.Range("D16").Value = "'0.00001339"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Value = "'1.665.94"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "'95.27"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'0.06997"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'7.012"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "'0.9982"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "'24.414.81"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "'3.126"
$ws.Range("E25").Value = "  +16.12%  "
$ws.Range("D26").Value = "'2.431"
$ws.Range("E26").Value = "  +4.73%  "
$ws.Range("D27").Value = "'22.57"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'157.88"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'142.89"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "'5.442"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("D31").Value = "'8.075"
$ws.Range("E31").Value = "  -6.85%  "
$ws.Range("D32").Value = "'2.536"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").Value = "'1.851.70"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "'1.068"
$ws.Range("E34").Value = "  +11.00%  "
$ws.Range("D35").Value = "'0.08277"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("D36").Value = "'0.03056"
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("D37").Value = "'6.921"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "'11.22"
$ws.Range("E38").Value = "  +12.63%  "
$ws.Range("D39").Value = "'0.2774"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("D40").Value = "'0.09258"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "'13.84"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("D42").Value = "'0.7724"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").Value = "'1.453"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "'16.65"
$ws.Range("E44").Value = "  +4.06%  "
$ws.Range("D45").Value = "'0.7136"
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").Value = "'2.547"
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("D47").Value = "'4.144"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "'1.0000"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'0.08436"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "'136.74"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "'1.272"
$ws.Range("E51").Value = "  +1.28%  "
